# OLX Monitor update — 2026-02-16 09:38 check run
$wb = $excel.ActiveWorkbook
$timestamp = "2026-02-16 09:38"

# ---------------------------------------------------------------------------
# PODSUMOWANIE (summary) sheet: refresh the "last checked" timestamp for every
# profile; wszystkie-lublin also picked up 2 new listings (432 -> 434).
# ---------------------------------------------------------------------------
$sum = $wb.Worksheets.Item("PODSUMOWANIE")
$sum.Range("B2").Value = $timestamp
$sum.Range("C2").Value = 434
$sum.Range("D2").Value = 2
$sum.Range("B3").Value = $timestamp
$sum.Range("B4").Value = $timestamp
$sum.Range("B5").Value = $timestamp
$sum.Range("B6").Value = $timestamp

# ---------------------------------------------------------------------------
# Per-profile detail sheets: append a new check row (row 14) underneath the
# existing history (rows 2-13). Row 14 is an "even" row, so its banding
# matches row 12's template rather than row 13's (odd-row) template.
# ---------------------------------------------------------------------------
function Add-CheckRow($Sheet, $NewRow, $TemplateEvenRow, $Date, $Total, $New, $Removed, $Net, $NewDetails, $RemovedDetails, $Status, $IValue) {
    $srcRange = $Sheet.Range("A" + $TemplateEvenRow + ":I" + $TemplateEvenRow)
    $dstRange = $Sheet.Range("A" + $NewRow + ":I" + $NewRow)
    $srcRange.Copy($dstRange)
    $Sheet.Rows.Item($NewRow).RowHeight = 18

    $Sheet.Range("A" + $NewRow).Value = $Date
    $Sheet.Range("B" + $NewRow).Value = $Total
    $Sheet.Range("C" + $NewRow).Value = $New
    $Sheet.Range("D" + $NewRow).Value = $Removed
    $Sheet.Range("E" + $NewRow).Value = $Net
    $Sheet.Range("F" + $NewRow).Value = $NewDetails
    $Sheet.Range("G" + $NewRow).Value = $RemovedDetails
    $Sheet.Range("H" + $NewRow).Value = $Status
    if ($IValue -ne "") {
        $Sheet.Range("I" + $NewRow).Value = $IValue
    }

    if ($New -gt 0) {
        $Sheet.Range("C4").Copy($Sheet.Range("C" + $NewRow))
        $Sheet.Range("C" + $NewRow).Value = $New
        $Sheet.Range("F4").Copy($Sheet.Range("F" + $NewRow))
        $Sheet.Range("F" + $NewRow).Value = $NewDetails
    }
}

$wszystkie = $wb.Worksheets.Item("wszystkie-lublin")
Add-CheckRow $wszystkie 14 12 $timestamp 434 2 0 2 "—" "—" "OK" ""

$artymiuk = $wb.Worksheets.Item("artymiuk")
Add-CheckRow $artymiuk 14 12 $timestamp 0 0 0 0 "—" "—" "OK" ""

$poqui = $wb.Worksheets.Item("poqui")
Add-CheckRow $poqui 14 12 $timestamp 5 0 0 0 "—" "—" "OK" "1951OR|183ger|17vbYq|17NeTz|18KAEc"

$stylowe = $wb.Worksheets.Item("stylowepokoje")
Add-CheckRow $stylowe 14 12 $timestamp 2 0 0 0 "—" "—" "OK" "195dLc|16ZeYm"

$villahome = $wb.Worksheets.Item("villahome")
Add-CheckRow $villahome 14 12 $timestamp 0 0 0 0 "—" "—" "OK" ""

Write-Output "OLX monitor updated: $timestamp"
